$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I18").Value = "[4]"
$ws.Range("I19").Value = "[4]"

$ws.Columns.Item(6).ColumnWidth = 54
$ws.Rows.Item(19).RowHeight = 13.8

[void]$ws.Range("I19").Select()
